$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Epoch Accuracy values in column B (re-run of training produced new numbers)
$ws.Range("B3").Value = 0.71875
$ws.Range("B4").Value = 0.671875
$ws.Range("B5").Value = 0.578125
$ws.Range("B6").Value = 0.515625
$ws.Range("B7").Value = 0.546875
$ws.Range("B8").Value = 0.53125
$ws.Range("B9").Value = 0.484375
$ws.Range("B10").Value = 0.515625
$ws.Range("B11").Value = 0.59375
$ws.Range("B12").Value = 0.515625
$ws.Range("B13").Value = 0.53125
$ws.Range("B14").Value = 0.5
$ws.Range("B17").Value = 0.5
$ws.Range("B18").Value = 0.5
$ws.Range("B19").Value = 0.5
$ws.Range("B21").Value = 0.5
$ws.Range("B22").Value = 0.5
$ws.Range("B23").Value = 0.5
$ws.Range("B24").Value = 0.5
$ws.Range("B25").Value = 0.5
$ws.Range("B26").Value = 0.5
$ws.Range("B27").Value = 0.5
$ws.Range("B28").Value = 0.5
$ws.Range("B29").Value = 0.515625
$ws.Range("B30").Value = 0.484375
$ws.Range("B31").Value = 0.515625
$ws.Range("B32").Value = 0.515625
$ws.Range("B33").Value = 0.515625
$ws.Range("B34").Value = 0.515625
$ws.Range("B35").Value = 0.515625
$ws.Range("B36").Value = 0.515625
$ws.Range("B37").Value = 0.515625
$ws.Range("B38").Value = 0.515625
$ws.Range("B39").Value = 0.515625
$ws.Range("B40").Value = 0.515625
$ws.Range("B41").Value = 0.515625
$ws.Range("B42").Value = 0.515625
$ws.Range("B43").Value = 0.515625
$ws.Range("B44").Value = 0.515625
$ws.Range("B45").Value = 0.515625
$ws.Range("B46").Value = 0.515625
$ws.Range("B47").Value = 0.515625
$ws.Range("B48").Value = 0.515625
$ws.Range("B49").Value = 0.515625
$ws.Range("B50").Value = 0.515625
$ws.Range("B51").Value = 0.515625
$ws.Range("B52").Value = 0.515625
$ws.Range("B53").Value = 0.515625
$ws.Range("B54").Value = 0.515625
$ws.Range("B55").Value = 0.515625
$ws.Range("B56").Value = 0.515625
$ws.Range("B57").Value = 0.515625
$ws.Range("B58").Value = 0.515625
$ws.Range("B59").Value = 0.515625
$ws.Range("B60").Value = 0.515625
$ws.Range("B61").Value = 0.515625
$ws.Range("B103").Value = 0.4375
$ws.Range("B104").Value = 0.40625
$ws.Range("B105").Value = 0.515625
$ws.Range("B106").Value = 0.34375
$ws.Range("B107").Value = 0.25
$ws.Range("B108").Value = 0.328125
$ws.Range("B109").Value = 0.515625
$ws.Range("B110").Value = 0.46875
$ws.Range("B111").Value = 0.421875
$ws.Range("B112").Value = 0.59375
$ws.Range("B113").Value = 0.5
$ws.Range("B114").Value = 0.375
$ws.Range("B115").Value = 0.4375
$ws.Range("B116").Value = 0.34375
$ws.Range("B117").Value = 0.40625
$ws.Range("B118").Value = 0.3934426229508197

# Refresh the <__main__.DisplayOutputs ...> repr text in column A (rows 102-118) to match the new run
$displayText = "<__main__.DisplayOutputs object at 0x7f612057f4c0>"
for ($row = 102; $row -le 118; $row++) {
    $ws.Cells.Item($row, 1).Value = $displayText
}

# Restore the worksheet selection state saved with the workbook
[void]$ws.Range("A2:B118").Select()
